$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 1989
$ws.Range("A4").Value = 1989
$ws.Range("A5").Value = 1989
$ws.Range("A6").Value = 1990

$ws.Range("A7").Select()
